$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet

# Insert a new worksheet right after "Arkusz1" and make it the active sheet
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Arkusz2"

# --- Header row (row 1) ---
$ws2.Range("A1").Value = "p1depth"
$ws2.Range("B1").Value = "n+"
$ws2.Range("C1").Value = "n-"
$ws2.Range("D1").Value = "mean-"
$ws2.Range("E1").Value = "mean+"
$ws2.Range("F1").Value = "(1+mean+ /mean- )"
$ws2.Range("G1").Value = "(1+mean- /mean+ )"
$ws2.Range("I1").Value = "p0"
$ws2.Range("J1").Value = "1-p0"
$ws2.Range("K1").Value = "p+"
$ws2.Range("L1").Value = "p-"
$ws2.Range("M1").Value = "Suma p0,p+,p-"

# --- Data / formula row (row 2) ---
$ws2.Range("A2").Value = 5
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Formula = "=2+4+8+31"
$ws2.Range("D2").Formula = "=(1*31+2*8+3*4+4*2)/C2"
$ws2.Range("E2").Value = 6
$ws2.Range("F2").Formula = "=(1+E2/D2)"
$ws2.Range("G2").Formula = "=(1+D2/E2)"
$ws2.Range("I2").Formula = "=1/A2"
$ws2.Range("J2").Formula = "=1-I2"
$ws2.Range("K2").Formula = "=J2/(B2*F2)"
$ws2.Range("L2").Formula = "=J2/G2"
$ws2.Range("M2").Formula = "=I2+K2+L2"

# --- Program listing (rows 6-13) ---
$ws2.Range("A6").Value = "Program:"
$ws2.Range("A7").Value = "1 -> 31 (16*1 + 15*'+')"
$ws2.Range("A8").Value = "2 -> 8"
$ws2.Range("A9").Value = "3 -> 4"
$ws2.Range("A10").Value = "4 -> 2"
$ws2.Range("A11").Value = "5 -> 1"
$ws2.Range("A12").Value = "6 -> write( *5* )"
$ws2.Range("A13").Value = "7 -> program"

# Match the selection left on the new sheet by the author
$ws2.Range("D17").Select() | Out-Null
